# Apply the "Added a few more slots" edit to the Gods of Giza review doc:
#  1. Remove the old "Meta description: ..." paragraph from the top of the
#     document (right under the H1 title).
#  2. Insert a new bold "Play Gods of Giza Slot for Free - Review & Demo"
#     paragraph right before the final (italic) image-prompt paragraph.
#  3. Replace the text of that final italic paragraph with the meta
#     description copy (keeping its italic run formatting intact).

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph -----------------------
# (the bold-lead-in paragraph right under the H1 title)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith("Meta description")) {
        $candidate.Range.Delete()
        break
    }
}

# --- Step 2: insert a new bold heading paragraph before the last paragraph
$count = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs.Item($count - 1)
$insertPoint = $beforeLast.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$newXml = "<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gods of Giza Slot for Free - Review &amp; Demo</w:t></w:r></w:p>"
[void]$newPara.Range.InsertXML($newXml)

# --- Step 3: swap the text of the (now) final paragraph --------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$startPos = $lastPara.Range.Start
$endPos = $lastPara.Range.End - 1
$textRange = $d.Range($startPos, $endPos)
$textRange.Text = "Try Gods of Giza slot game for free with our review. Features, graphics, and bonus options of Pragmatic Play's Gods of Giza Slot explained."
